$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Mark the to-do item in row 11 ("Intro: change section 1.3 ...") as
# completed -> the D column checkbox flips from FALSE to TRUE, which in
# turn recalculates the per-row F column and the overall completion
# percentage in row 57 (D57/F57).
$ws.Range("D11").Value = $true

# Update the saved view state: scroll position resets to the top of the
# sheet and the active selection moves to C2.
$ws.Range("C2").Select() | Out-Null

$wb.Save()
